$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value, taken from the Tue Jul 18 22:11:05 UTC 2023
# GitHub Actions refresh of the cryptos list.
$updates = @{
    'D2' = '29.831.75'
    'E2' = '  -1.08%  '
    'D3' = '1.899.52'
    'E4' = '  +0.19%  '
    'D5' = '0.7676'
    'E5' = '  +3.81%  '
    'D6' = '240.59'
    'E6' = '  -1.33%  '
    'E7' = '  +0.19%  '
    'D8' = '0.3053'
    'E8' = '  -2.14%  '
    'D9' = '25.36'
    'E9' = '  -4.57%  '
    'D10' = '0.06850'
    'E10' = '  -1.59%  '
    'D11' = '0.07986'
    'E11' = '  +0.23%  '
    'D12' = '1.892.26'
    'E12' = '  -1.13%  '
    'D13' = '0.7362'
    'E13' = '  -5.51%  '
    'D14' = '5.171'
    'E14' = '  -1.95%  '
    'D15' = '91.17'
    'E15' = '  -1.15%  '
    'D16' = '29.836.44'
    'E16' = '  -1.12%  '
    'D17' = '13.78'
    'E17' = '  -3.81%  '
    'D18' = '5.889'
    'E18' = '  +0.13%  '
    'D19' = '245.53'
    'E19' = '  +1.34%  '
    'D20' = '0.000007696'
    'E20' = '  -1.59%  '
    'E21' = '  +0.21%  '
    'D22' = '2.132.72'
    'E22' = '  -0.33%  '
    'D23' = '1.001'
    'E23' = '  +0.17%  '
    'D24' = '6.884'
    'E24' = '  -1.99%  '
    'D25' = '167.21'
    'E25' = '  -0.56%  '
    'D26' = '9.247'
    'E26' = '  -1.52%  '
    'D27' = '18.68'
    'E27' = '  -2.00%  '
    'D28' = '0.1287'
    'E28' = '  +0.35%  '
    'D29' = '2.030'
    'E29' = '  -1.71%  '
    'D30' = '1.399'
    'E30' = '  +3.81%  '
    'D31' = '1.513'
    'E31' = '  -1.93%  '
    'D32' = '4.274'
    'E32' = '  -1.37%  '
    'E33' = '  -0.80%  '
    'D34' = '0.05279'
    'E34' = '  +2.53%  '
    'E35' = '  -4.26%  '
    'D36' = '0.7255'
    'E36' = '  -2.61%  '
    'E37' = '  -0.14%  '
    'E38' = '  -1.50%  '
    'D39' = '2.778'
    'E39' = '  -0.85%  '
    'D40' = '6.207'
    'E40' = '  -2.29%  '
    'D41' = '0.4408'
    'E41' = '  -1.85%  '
    'D42' = '72.04'
    'E42' = '  -4.02%  '
    'E43' = '  +0.09%  '
    'D44' = '0.8334'
    'E44' = '  -0.53%  '
    'D45' = '1.882'
    'E45' = '  -3.70%  '
    'D46' = '7.600'
    'E46' = '  -2.79%  '
    'D47' = '99.97'
    'E47' = '  -1.39%  '
    'D48' = '9.729'
    'E48' = '  -1.84%  '
    'D49' = '2.040.99'
    'E49' = '  -0.63%  '
    'D50' = '36.17'
    'E50' = '  -2.76%  '
    'B51' = 'NEARProtocol'
    'C51' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D51' = '1.471'
    'E51' = '  -1.25%  '
}

# Cells in the Price column (D) hold numeric-looking strings that must
# stay literal text (matching the source inlineStr cells), e.g. trailing
# zeros like '0.06850' which Excel would otherwise normalize to 0.0685
# if assigned as a plain numeric-looking value. Prefixing with a leading
# apostrophe forces text entry; re-applying the 'Normal' style afterwards
# clears the quotePrefix flag Excel stamps on the cell format so no stray
# style diff is left behind (source cells carry no explicit style).
foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    $range = $ws.Range($cell)
    if ($cell.Substring(0,1) -eq 'D') {
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
